{"js": "// Add a \"Meta description\" paragraph right after the title (Heading1) and\n// replace the old \"Play Age of Asgard...\" + italic blurb pairing at the\n// bottom of the document with just an updated italic \"Prompt: ...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) Insert the new \"Meta description: ...\" paragraph right after the\n//        document title (first paragraph). ---\nconst titlePara = paragraphs.items[0];\nconst metaLabel = \"Meta description\";\nconst metaRest = \": Read our review of Age of Asgard slot game. Play for free and enjoy special features triggered by symbol interaction, unique gameplay and Norse theme.\";\n\nconst metaPara = titlePara.insertParagraph(metaLabel, \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Append the (non-bold) remainder of the sentence.\nmetaPara.insertText(metaRest, \"End\");\nawait context.sync();\n\n// Bold only the \"Meta description\" label, matching the source formatting.\nconst boldHits = metaPara.search(metaLabel, { matchCase: true });\nboldHits.load(\"items\");\nawait context.sync();\nboldHits.items[0].font.bold = true;\nawait context.sync();\n\n// --- 2) Drop the stray duplicate title paragraph and refresh the prompt\n//        text of the italic blurb near the end of the document. ---\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nconst count = allParagraphs.items.length;\nconst duplicateTitlePara = allParagraphs.items[count - 2];\nconst blurbPara = allParagraphs.items[count - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst newBlurb = \"Prompt: Create a feature image for Age of Asgard, a slot game that offers a unique twist on the beloved mythological theme of Norse gods and their battles. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The image should showcase the excitement and epicness of the game, with an ancient Viking village in the background and flames and warriors on both sides of the grid. The Maya warrior should be dressed in a traditional Viking outfit and have a big smile on their face, holding up a winning combination on the slot machine. Be sure to incorporate elements of Norse mythology and the different symbols of the game into the design of the image.\";\n\nblurbPara.getRange().insertText(newBlurb, \"Replace\");\nawait context.sync();\n", "ps1": "# Add a \"Meta description\" paragraph right after the title (Heading1) and\n# replace the old \"Play Age of Asgard...\" + italic blurb pairing at the\n# bottom of the document with just an updated italic \"Prompt: ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert the new \"Meta description: ...\" paragraph right after the\n#        document title (first paragraph). ---\n$metaLabel = \"Meta description\"\n$metaRest = \": Read our review of Age of Asgard slot game. Play for free and enjoy special features triggered by symbol interaction, unique gameplay and Norse theme.\"\n\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n$metaPara.Range.InsertAfter($metaLabel)\n$metaPara.Range.InsertAfter($metaRest)\n\n$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $metaLabel.Length)\n$boldRange.Bold = 1\n\n# --- 2) Drop the stray duplicate title paragraph and refresh the prompt\n#        text of the italic blurb near the end of the document. ---\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs.Item($count - 1)\n$blurbPara = $d.Paragraphs.Item($count)\n\n$newBlurb = \"Prompt: Create a feature image for Age of Asgard, a slot game that offers a unique twist on the beloved mythological theme of Norse gods and their battles. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The image should showcase the excitement and epicness of the game, with an ancient Viking village in the background and flames and warriors on both sides of the grid. The Maya warrior should be dressed in a traditional Viking outfit and have a big smile on their face, holding up a winning combination on the slot machine. Be sure to incorporate elements of Norse mythology and the different symbols of the game into the design of the image.\"\n\n$blurbPara.Range.Text = $newBlurb\n$duplicateTitlePara.Range.Delete()\n"}
